# Refresh the cryptocurrency price / volume(1h) snapshot (GitHub Actions scrape).
# Rows 28/29 (PEPE <-> PancakeSwap) also swapped ranking order this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.245.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.265.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "496.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.337"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.667.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.221.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.271.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0691"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.943"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.57%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.374"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "125.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0495"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.546"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "241.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.373"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -0.68%  "
